$d = $word.ActiveDocument

$replacements = @(
    @("54÷6=", "96÷4="),
    @("47÷7=", "49÷9="),
    @("94÷3=", "19÷6="),
    @("94÷2=", "98÷5="),
    @("82÷6=", "88÷9="),
    @("39÷5=", "43÷8="),
    @("56÷5=", "99÷5="),
    @("54÷9=", "15÷6="),
    @("48÷5=", "23÷9="),
    @("74÷8=", "56÷6="),
    @("45÷8=", "55÷4="),
    @("68÷5=", "80÷8="),
    @("59÷2=", "79÷8="),
    @("43÷2=", "28÷5="),
    @("58÷6=", "48÷4="),
    @("88÷2=", "25÷9="),
    @("62÷4=", "32÷8="),
    @("75÷6=", "61÷5="),
    @("22÷4=", "30÷5="),
    @("52÷2=", "92÷9="),
    @("11÷5=", "55÷6="),
    @("26÷3=", "86÷3="),
    @("53÷3=", "40÷6="),
    @("73÷7=", "44÷7="),
    @("81÷8=", "39÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
